$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,36
$row2[0,0] = 47.41391982
$row2[0,1] = 49.77440937
$row2[0,2] = 53.27488453
$row2[0,3] = 57.03343688919741
$row2[0,4] = 61.13156901037621
$row2[0,5] = 64.43568065268151
$row2[0,6] = 67.54041373905439
$row2[0,7] = 70.78205234958706
$row2[0,8] = 74.53350127977443
$row2[0,9] = 78.63284378575129
$row2[0,10] = 82.95765011613798
$row2[0,11] = 87.83141238184376
$row2[0,12] = 92.99150767678005
$row2[0,13] = 98.45475847099397
$row2[0,14] = 104.2389759390994
$row2[0,15] = 110.3630157829019
$row2[0,16] = 116.8468426421195
$row2[0,17] = 124.1497704414409
$row2[0,18] = 131.9091308592002
$row2[0,19] = 140.1534519404672
$row2[0,20] = 148.9130426867465
$row2[0,21] = 158.2201079217626
$row2[0,22] = 168.1088642307585
$row2[0,23] = 178.615668245181
$row2[0,24] = 189.7791476111465
$row2[0,25] = 201.6403443032959
$row2[0,26] = 214.2428659899882
$row2[0,27] = 227.6330449130789
$row2[0,28] = 241.8601107233551
$row2[0,29] = 256.9763676100175
$row2[0,30] = 273.0373902166239
$row2[0,31] = 290.1022270716156
$row2[0,32] = 308.2336163642333
$row2[0,33] = 327.4982172863562
$row2[0,34] = 347.9668558667535
$row2[0,35] = 369.7147847609926
$ws.Range("J2:AS2").Value = $row2

$row3 = New-Object 'object[,]' 1,36
$row3[0,0] = 86610000
$row3[0,1] = 86610000
$row3[0,2] = 86610000
$row3[0,3] = 86610000
$row3[0,4] = 86610000
$row3[0,5] = 86610000
$row3[0,6] = 86610000
$row3[0,7] = 86610000
$row3[0,8] = 86610000
$row3[0,9] = 86610000
$row3[0,10] = 86610000
$row3[0,11] = 86610000
$row3[0,12] = 86610000
$row3[0,13] = 86610000
$row3[0,14] = 86610000
$row3[0,15] = 86610000
$row3[0,16] = 86610000
$row3[0,17] = 86610000
$row3[0,18] = 86610000
$row3[0,19] = 86610000
$row3[0,20] = 86610000
$row3[0,21] = 86610000
$row3[0,22] = 86610000
$row3[0,23] = 86610000
$row3[0,24] = 86610000
$row3[0,25] = 86610000
$row3[0,26] = 86610000
$row3[0,27] = 86610000
$row3[0,28] = 86610000
$row3[0,29] = 86610000
$row3[0,30] = 86610000
$row3[0,31] = 86610000
$row3[0,32] = 86610000
$row3[0,33] = 86610000
$row3[0,34] = 86610000
$row3[0,35] = 86610000
$ws.Range("J3:AS3").Value = $row3

$row10 = New-Object 'object[,]' 1,36
$row10[0,0] = 3.919
$row10[0,1] = 3.919
$row10[0,2] = 3.919
$row10[0,3] = 3.919
$row10[0,4] = 3.919
$row10[0,5] = 3.919
$row10[0,6] = 3.919
$row10[0,7] = 3.919
$row10[0,8] = 3.919
$row10[0,9] = 3.919
$row10[0,10] = 3.919
$row10[0,11] = 3.919
$row10[0,12] = 3.919
$row10[0,13] = 3.919
$row10[0,14] = 3.919
$row10[0,15] = 3.919
$row10[0,16] = 3.919
$row10[0,17] = 3.919
$row10[0,18] = 3.919
$row10[0,19] = 3.919
$row10[0,20] = 3.919
$row10[0,21] = 3.919
$row10[0,22] = 3.919
$row10[0,23] = 3.919
$row10[0,24] = 3.919
$row10[0,25] = 3.919
$row10[0,26] = 3.919
$row10[0,27] = 3.919
$row10[0,28] = 3.919
$row10[0,29] = 3.919
$row10[0,30] = 3.919
$row10[0,31] = 3.919
$row10[0,32] = 3.919
$row10[0,33] = 3.919
$row10[0,34] = 3.919
$row10[0,35] = 3.919
$ws.Range("J10:AS10").Value = $row10

$row11 = New-Object 'object[,]' 1,36
$row11[0,0] = 35930358.65
$row11[0,1] = 36812067.36
$row11[0,2] = 37669089.91
$row11[0,3] = 38469814.97
$row11[0,4] = 39216539.25
$row11[0,5] = 39967867.44
$row11[0,6] = 40725784.39
$row11[0,7] = 41471864.08
$row11[0,8] = 42210184.93
$row11[0,9] = 42945415.32
$row11[0,10] = 43669534.48
$row11[0,11] = 44382344.6
$row11[0,12] = 45084671.91
$row11[0,13] = 45777069.76
$row11[0,14] = 46458023.83
$row11[0,15] = 47127449.01
$row11[0,16] = 47788606.49
$row11[0,17] = 48441965.2
$row11[0,18] = 49086670.5
$row11[0,19] = 49727259.72
$row11[0,20] = 50359967.09
$row11[0,21] = 50982558.35
$row11[0,22] = 51597751.77
$row11[0,23] = 52206777.06
$row11[0,24] = 52812186.15
$row11[0,25] = 53394079.62
$row11[0,26] = 53950829.14
$row11[0,27] = 54486246.79
$row11[0,28] = 55001544.1
$row11[0,29] = 55494147.17
$row11[0,30] = 55959041.35
$row11[0,31] = 56402164.38
$row11[0,32] = 56817570.33
$row11[0,33] = 57205024.46
$row11[0,34] = 57567278.61
$row11[0,35] = 57906580.52
$ws.Range("J11:AS11").Value = $row11

$row12 = New-Object 'object[,]' 1,36
$row12[0,0] = 16612464.35
$row12[0,1] = 17589734.64
$row12[0,2] = 18597942.09
$row12[0,3] = 19620628.03
$row12[0,4] = 20656039.76
$row12[0,5] = 21736650.56
$row12[0,6] = 22862549.61
$row12[0,7] = 24025883.92
$row12[0,8] = 25227921.07
$row12[0,9] = 26473657.68
$row12[0,10] = 27758275.52
$row12[0,11] = 29080114.4
$row12[0,12] = 30441632.09
$row12[0,13] = 31842758.24
$row12[0,14] = 33281813.17
$row12[0,15] = 34757854.99
$row12[0,16] = 36273745.51
$row12[0,17] = 37828920.8
$row12[0,18] = 39423166.5
$row12[0,19] = 41057331.28
$row12[0,20] = 42728568.91
$row12[0,21] = 44431433.65
$row12[0,22] = 46167850.23
$row12[0,23] = 47934739.94
$row12[0,24] = 49727783.85
$row12[0,25] = 51563617.38
$row12[0,26] = 53437512.86
$row12[0,27] = 55349556.21
$row12[0,28] = 57306096.9
$row12[0,29] = 59300745.83
$row12[0,30] = 61330921.65
$row12[0,31] = 63398524.62
$row12[0,32] = 65502279.68
$row12[0,33] = 67639532.54000001
$row12[0,34] = 69813684.39
$row12[0,35] = 72024939.48
$ws.Range("J12:AS12").Value = $row12
